{"js": "// Replace the static \"KABUPATEN HULU SUNGAI TENGAH\" heading with the\n// templated placeholder \"${u_kabupaten}\" (keeps the existing bold/italic\n// Arial run formatting because insertText(\"Replace\") reuses the matched\n// range's formatting).\nconst heading = context.document.body.search(\"KABUPATEN HULU SUNGAI TENGAH\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nheading.load(\"items\");\nawait context.sync();\n\nif (heading.items.length > 0) {\n  heading.items[0].insertText(\"${u_kabupaten}\", \"Replace\");\n  await context.sync();\n}\n\n// Replace the hard-coded \"Kabupaten Hulu Sungai Tengah\" in the\n// \"Kepala BPS Kabupaten Hulu Sungai Tengah\" signature line with the\n// templated placeholder \"${kabupaten}\".\nconst signature = context.document.body.search(\" BPS Kabupaten Hulu Sungai Tengah\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nsignature.load(\"items\");\nawait context.sync();\n\nif (signature.items.length > 0) {\n  signature.items[0].insertText(\" BPS ${kabupaten}\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Replace the static \"KABUPATEN HULU SUNGAI TENGAH\" heading with the\n# templated placeholder \"${u_kabupaten}\" and the hard-coded \"Kabupaten Hulu\n# Sungai Tengah\" in the \"Kepala BPS ...\" signature line with \"${kabupaten}\".\n#\n# NOTE: replacement strings use single quotes throughout so PowerShell does\n# not try to interpolate the literal \"${...}\" template syntax as a variable.\n\n$d = $word.ActiveDocument\n\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Text = 'KABUPATEN HULU SUNGAI TENGAH'\n$find1.Replacement.ClearFormatting()\n$find1.Replacement.Text = '${u_kabupaten}'\n$find1.Execute($find1.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2)\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = ' BPS Kabupaten Hulu Sungai Tengah'\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = ' BPS ${kabupaten}'\n$find2.Execute($find2.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n"}
